$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-CellText 2 4 "65.577.92"
Set-CellText 2 5 "  +1.92%  "
Set-CellText 3 4 "2.644.93"
Set-CellText 3 5 "  +0.41%  "
Set-CellText 4 5 "  -0.08%  "
Set-CellText 5 4 "603.75"
Set-CellText 5 5 "  +1.13%  "
Set-CellText 6 4 "156.39"
Set-CellText 6 5 "  +2.65%  "
Set-CellText 7 5 "  -0.09%  "
Set-CellText 8 5 "  -0.49%  "
Set-CellText 9 4 "2.644.01"
Set-CellText 9 5 "  +0.51%  "
Set-CellText 10 4 "0.124"
Set-CellText 10 5 "  +7.69%  "
Set-CellText 11 5 "  +1.77%  "
Set-CellText 12 4 "5.86"
Set-CellText 12 5 "  +0.40%  "
Set-CellText 13 5 "  +1.45%  "
Set-CellText 14 4 "29.67"
Set-CellText 14 5 "  +5.62%  "
Set-CellText 15 4 "0.0000195"
Set-CellText 15 5 "  +13.58%  "
Set-CellText 16 4 "3.121.62"
Set-CellText 16 5 "  +0.39%  "
Set-CellText 17 4 "65.319.09"
Set-CellText 17 5 "  +1.68%  "
Set-CellText 18 4 "2.652.42"
Set-CellText 18 5 "  -0.26%  "
Set-CellText 19 4 "12.60"
Set-CellText 19 5 "  +2.30%  "
Set-CellText 20 4 "4.86"
Set-CellText 20 5 "  +1.52%  "
Set-CellText 21 4 "357.15"
Set-CellText 21 5 "  +1.79%  "
Set-CellText 22 4 "7.41"
Set-CellText 22 5 "  +4.29%  "
Set-CellText 23 5 "  +0.00%  "
Set-CellText 24 4 "69.45"
Set-CellText 24 5 "  +2.52%  "
Set-CellText 25 5 "  +0.51%  "
Set-CellText 26 2 "InternetComputer(DFINITY)"
Set-CellText 26 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-CellText 26 4 "9.39"
Set-CellText 26 5 "  +1.69%  "
Set-CellText 27 2 "PEPE"
Set-CellText 27 3 "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-CellText 27 4 "0.0000105"
Set-CellText 27 5 "  +15.20%  "
Set-CellText 28 5 "  -3.20%  "
Set-CellText 29 5 "  +1.72%  "
Set-CellText 30 4 "8.11"
Set-CellText 30 5 "  -2.74%  "
Set-CellText 31 5 "  +0.12%  "
Set-CellText 32 5 "  +4.50%  "
Set-CellText 33 4 "529.26"
Set-CellText 33 5 "  -4.42%  "
Set-CellText 34 4 "1.78"
Set-CellText 34 5 "  -3.45%  "
Set-CellText 35 4 "5.53"
Set-CellText 35 5 "  -0.08%  "
Set-CellText 36 4 "6.34"
Set-CellText 36 5 "  +2.07%  "
Set-CellText 37 5 "  +1.80%  "
Set-CellText 38 4 "20.63"
Set-CellText 38 5 "  +2.77%  "
Set-CellText 39 4 "161.71"
Set-CellText 39 5 "  -2.65%  "
Set-CellText 40 2 "Stacks"
Set-CellText 40 3 "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-CellText 40 4 "1.98"
Set-CellText 40 5 "  -0.92%  "
Set-CellText 41 2 "FirstDigitalUSD"
Set-CellText 41 3 "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-CellText 41 4 "0.999"
Set-CellText 41 5 "  -0.08%  "
Set-CellText 42 5 "  -0.02%  "
Set-CellText 43 4 "41.96"
Set-CellText 43 5 "  +4.05%  "
Set-CellText 44 4 "164.53"
Set-CellText 44 5 "  -2.70%  "
Set-CellText 45 5 "  +0.45%  "
Set-CellText 46 4 "2.36"
Set-CellText 46 5 "  +6.51%  "
Set-CellText 47 4 "0.0607"
Set-CellText 47 5 "  +3.03%  "
Set-CellText 48 4 "22.81"
Set-CellText 48 5 "  -1.78%  "
Set-CellText 49 4 "0.651"
Set-CellText 49 5 "  +1.70%  "
Set-CellText 50 5 "  +3.27%  "
Set-CellText 51 5 "  +0.57%  "
